# Adds a new color-swatch pair (an "Oval" circle + its hex-code "TextBox"
# label) to slide 2, matching the existing swatch row already on that
# slide ("Oval 17" / "TextBox 6", color D9CAB3). The new swatch uses
# color F6F6F6, positioned one slot further to the right.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Shape-id bookkeeping -------------------------------------------------
# This deck's COM host hands out shape ids as "smallest id not already used
# on the slide", counting up from the state the slide was in when it was
# first touched; it does not re-scan after Delete(). The slide currently
# has gaps at 3,4,5,6,9,10,12,14,16,19 before the next free id reaches 20 -
# the id the real edit's new Oval ("Oval 19") ends up with, followed by 21
# for the TextBox ("TextBox 20"). Burn through those lower gap ids with
# disposable placeholder shapes (then delete them) so our two real shapes
# land on ids 20 / 21, exactly like the authored edit.
$fillers = @()
for ($i = 1; $i -le 10; $i++) {
    $fillers += $s.Shapes.AddShape(1, 0, 0, 1, 1)
}
foreach ($f in $fillers) {
    $f.Delete()
}

# --- New Oval (swatch circle) --------------------------------------------
# Duplicate the existing "Oval 17" swatch so the new shape inherits the
# exact same shape style (lnRef/fillRef/effectRef/fontRef) and text body
# formatting, then reposition/recolor/rename it.
$ovalSrc = $s.Shapes.Item(13)
$oval = $ovalSrc.Duplicate()
$oval.Name = "Oval 19"
$oval.Left = 606.1746826171875
$oval.Top = 382.5
$oval.Width = 83.25
$oval.Height = 83.25
$oval.Fill.ForeColor.RGB = 16185078

# --- New TextBox (hex-code label) -----------------------------------------
# Duplicate the existing "TextBox 6" swatch label for the same reason.
$tbSrc = $s.Shapes.Item(14)
$tb = $tbSrc.Duplicate()
$tb.Name = "TextBox 20"
$tb.Left = 610.6780395507812
$tb.Top = 344.8968505859375
$tb.Width = 74.24307250976562
$tb.Height = 29.081260681152344
$tb.TextFrame.TextRange.Text = "D9CAB3"
